{"js": "// Replace each two-digit multiplication expression in the table cells\n// with its updated version, following the authoring diff exactly.\nconst replacements = [\n  [\"29\u00d760=1740\", \"95\u00d798=9310\"],\n  [\"91\u00d770=6370\", \"89\u00d783=7387\"],\n  [\"12\u00d786=1032\", \"50\u00d720=1000\"],\n  [\"13\u00d754=702\", \"97\u00d793=9021\"],\n  [\"21\u00d790=1890\", \"55\u00d760=3300\"],\n  [\"99\u00d773=7227\", \"96\u00d799=9504\"],\n  [\"13\u00d747=611\", \"36\u00d797=3492\"],\n  [\"35\u00d726=910\", \"68\u00d768=4624\"],\n  [\"75\u00d741=3075\", \"72\u00d743=3096\"],\n  [\"58\u00d737=2146\", \"65\u00d770=4550\"],\n  [\"15\u00d779=1185\", \"96\u00d743=4128\"],\n  [\"52\u00d757=2964\", \"45\u00d764=2880\"],\n  [\"65\u00d760=3900\", \"50\u00d737=1850\"],\n  [\"20\u00d799=1980\", \"95\u00d733=3135\"],\n  [\"78\u00d788=6864\", \"39\u00d771=2769\"],\n  [\"95\u00d757=5415\", \"49\u00d795=4655\"],\n  [\"12\u00d715=180\", \"86\u00d740=3440\"],\n  [\"55\u00d790=4950\", \"38\u00d726=988\"],\n  [\"47\u00d743=2021\", \"41\u00d778=3198\"],\n  [\"33\u00d774=2442\", \"99\u00d757=5643\"],\n  [\"82\u00d789=7298\", \"45\u00d791=4095\"],\n  [\"43\u00d735=1505\", \"18\u00d731=558\"],\n  [\"31\u00d765=2015\", \"14\u00d761=854\"],\n  [\"78\u00d749=3822\", \"93\u00d769=6417\"],\n  [\"29\u00d764=1856\", \"54\u00d789=4806\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n\n", "ps1": "# Update each two-digit multiplication expression in the table cells\n# to match the regenerated answer key, per the authoring diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"29\u00d760=1740\", \"95\u00d798=9310\")\n    ,@(\"91\u00d770=6370\", \"89\u00d783=7387\")\n    ,@(\"12\u00d786=1032\", \"50\u00d720=1000\")\n    ,@(\"13\u00d754=702\", \"97\u00d793=9021\")\n    ,@(\"21\u00d790=1890\", \"55\u00d760=3300\")\n    ,@(\"99\u00d773=7227\", \"96\u00d799=9504\")\n    ,@(\"13\u00d747=611\", \"36\u00d797=3492\")\n    ,@(\"35\u00d726=910\", \"68\u00d768=4624\")\n    ,@(\"75\u00d741=3075\", \"72\u00d743=3096\")\n    ,@(\"58\u00d737=2146\", \"65\u00d770=4550\")\n    ,@(\"15\u00d779=1185\", \"96\u00d743=4128\")\n    ,@(\"52\u00d757=2964\", \"45\u00d764=2880\")\n    ,@(\"65\u00d760=3900\", \"50\u00d737=1850\")\n    ,@(\"20\u00d799=1980\", \"95\u00d733=3135\")\n    ,@(\"78\u00d788=6864\", \"39\u00d771=2769\")\n    ,@(\"95\u00d757=5415\", \"49\u00d795=4655\")\n    ,@(\"12\u00d715=180\", \"86\u00d740=3440\")\n    ,@(\"55\u00d790=4950\", \"38\u00d726=988\")\n    ,@(\"47\u00d743=2021\", \"41\u00d778=3198\")\n    ,@(\"33\u00d774=2442\", \"99\u00d757=5643\")\n    ,@(\"82\u00d789=7298\", \"45\u00d791=4095\")\n    ,@(\"43\u00d735=1505\", \"18\u00d731=558\")\n    ,@(\"31\u00d765=2015\", \"14\u00d761=854\")\n    ,@(\"78\u00d749=3822\", \"93\u00d769=6417\")\n    ,@(\"29\u00d764=1856\", \"54\u00d789=4806\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n\n"}
